$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D that need to hold number-like text (e.g. "582.81") must be
# forced to Text format, otherwise Excel auto-converts them to real numbers.
# Save original styles, switch to Text format, assign, then restore original style
# so no stray style/number-format change is left on the cell.
$numericTextCells = @("D5", "D6", "D7", "D13", "D18", "D21", "D23", "D24", "D26", "D27", "D28", "D32", "D34", "D37", "D38", "D40", "D42", "D43", "D44", "D45", "D47", "D48", "D50")
$origStyles = @{}
foreach ($addr in $numericTextCells) {
    $origStyles[$addr] = $ws.Range($addr).Style
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.576.08"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "2.637.31"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "582.81"
$ws.Range("E5").Value = "  -2.19%  "
$ws.Range("D6").Value = "157.11"
$ws.Range("E6").Value = "  +1.01%  "
$ws.Range("D7").Value = "0.645"
$ws.Range("E7").Value = "  +2.56%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("E10").Value = "  +1.17%  "
$ws.Range("D13").Value = "28.73"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("E14").Value = "  -4.24%  "
$ws.Range("D15").Value = "3.113.32"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "64.357.29"
$ws.Range("E16").Value = "  -1.58%  "
$ws.Range("D17").Value = "2.632.23"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "12.25"
$ws.Range("E18").Value = "  -2.48%  "
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").Value = "348.06"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "68.28"
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("D24").Value = "1.78"
$ws.Range("E24").Value = "  +8.64%  "
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").Value = "9.46"
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("D27").Value = "595.47"
$ws.Range("E27").Value = "  +10.57%  "
$ws.Range("D28").Value = "1.59"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").Value = "2.08"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").Value = "6.66"
$ws.Range("E34").Value = "  +4.65%  "
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("D37").Value = "20.07"
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("D40").Value = "152.74"
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "158.78"
$ws.Range("E42").Value = "  -0.90%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "2.41"
$ws.Range("E43").Value = "  +5.01%  "
$ws.Range("D44").Value = "4.02"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").Value = "23.49"
$ws.Range("E45").Value = "  +4.62%  "
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("D47").Value = "0.636"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").Value = "0.0256"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("E49").Value = "  +2.48%  "
$ws.Range("D50").Value = "19.26"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").Value = "0.0₆0237"
$ws.Range("E51").Value = "  -5.92%  "

# restore original styles on the cells we temporarily reformatted
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).Style = $origStyles[$addr]
}
